{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"2024-06-04 Tuesday\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"2024-06-05 Wednesday\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"702\u00f78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"390\u00f73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"266\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"625\u00f73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"588\u00f73=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"925\u00f74=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"385\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"143\u00f72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"987\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"509\u00f73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"641\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"741\u00f73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"383\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"724\u00f77=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"196\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"885\u00f77=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"498\u00f74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"925\u00f77=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"950\u00f77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"823\u00f75=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"922\u00f73=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"371\u00f79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"415\u00f77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"471\u00f79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"931\u00f75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"139\u00f72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"163\u00f74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"484\u00f73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"180\u00f75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"228\u00f79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"135\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"796\u00f75=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"205\u00f74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"893\u00f73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"233\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"123\u00f74=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"205\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"743\u00f76=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"195\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"398\u00f72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"592\u00f73=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"106\u00f72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"921\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"562\u00f79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"907\u00f74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"624\u00f72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"797\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"367\u00f75=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n{\n  const results = body.search(\"791\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"497\u00f78=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $oldText, $newText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nReplace-Text $d \"2024-06-04 Tuesday\" \"2024-06-05 Wednesday\"\nReplace-Text $d \"702\u00f78=\" \"390\u00f73=\"\nReplace-Text $d \"266\u00f72=\" \"625\u00f73=\"\nReplace-Text $d \"588\u00f73=\" \"925\u00f74=\"\nReplace-Text $d \"385\u00f72=\" \"143\u00f72=\"\nReplace-Text $d \"987\u00f76=\" \"509\u00f73=\"\nReplace-Text $d \"641\u00f79=\" \"741\u00f73=\"\nReplace-Text $d \"383\u00f79=\" \"724\u00f77=\"\nReplace-Text $d \"196\u00f79=\" \"885\u00f77=\"\nReplace-Text $d \"498\u00f74=\" \"925\u00f77=\"\nReplace-Text $d \"950\u00f77=\" \"823\u00f75=\"\nReplace-Text $d \"922\u00f73=\" \"371\u00f79=\"\nReplace-Text $d \"415\u00f77=\" \"471\u00f79=\"\nReplace-Text $d \"931\u00f75=\" \"139\u00f72=\"\nReplace-Text $d \"163\u00f74=\" \"484\u00f73=\"\nReplace-Text $d \"180\u00f75=\" \"228\u00f79=\"\nReplace-Text $d \"135\u00f72=\" \"796\u00f75=\"\nReplace-Text $d \"205\u00f74=\" \"893\u00f73=\"\nReplace-Text $d \"233\u00f76=\" \"123\u00f74=\"\nReplace-Text $d \"205\u00f79=\" \"743\u00f76=\"\nReplace-Text $d \"195\u00f79=\" \"398\u00f72=\"\nReplace-Text $d \"592\u00f73=\" \"106\u00f72=\"\nReplace-Text $d \"921\u00f76=\" \"562\u00f79=\"\nReplace-Text $d \"907\u00f74=\" \"624\u00f72=\"\nReplace-Text $d \"797\u00f76=\" \"367\u00f75=\"\nReplace-Text $d \"791\u00f79=\" \"497\u00f78=\"\n"}
